# Reproduce the authoring change: column D (the duplicate text column that
# shadowed A, with header "Mtetric") was removed from the "Test" worksheet.
# Deleting the column shifts E:H left into D:G and keeps all the other data
# in place, exactly matching the diff (dimension A1:H5 -> A1:G5, header and
# value cells shifted one column left, shared string "Mtetric" falling out
# of use, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test")

$ws.Columns.Item(4).Delete()

# Rows 4 and 5 previously only carried a text label in (old) column D and had
# no data at all under the metric columns. After the shift those rows need
# real numeric values filled into D:G (Change Fail % / 0-15% metric rows).
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0.5

$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0

# Update the chart's series so it points at the new D:G location instead of
# the old E:H location, and drop the now-orphaned series-name reference
# (the cell that used to hold "Lead Time for Changes (Median)" in $D$3 no
# longer exists after the column shift).
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Name = ""
$ser.Formula = "=SERIES(,Test!`$D`$1:`$G`$1,Test!`$D`$3:`$G`$3,1)"

# The chart title was removed from the chart.
$chart.HasTitle = $false

# The plot area now uses a manual (fixed) inner layout instead of the
# automatic one.
$pa = $chart.PlotArea()
$pa.InsideLeft = 0.11208055723803756
$pa.InsideTop = 0.07407407407407407
$pa.InsideWidth = 0.88791944276196244
$pa.InsideHeight = 0.8416746864975212

# The chart's right anchor moved in one column (from column J to column I)
# because column D was removed upstream of it; the top-left anchor point is
# unaffected since it sits before the deleted column. Resize the chart's
# width so it lands back on the same to-cell/offset the original file had.
$co.Width = 363.009765625

# The active selection left on the sheet after the edit.
[void]$ws.Range("J6").Select()
